$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; existing rows 27..106 shift down to 28..107
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with its data
$ws.Cells.Item(27, 1).Value  = 11
$ws.Cells.Item(27, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(27, 3).Value  = "Bíobío"
$ws.Cells.Item(27, 4).Value  = 44623
$ws.Cells.Item(27, 5).Value  = 8
$ws.Cells.Item(27, 6).Value  = 100112043
$ws.Cells.Item(27, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(27, 8).Value  = "Sin especificar"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 220
$ws.Cells.Item(27, 11).Value = 14000
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = 14455
$ws.Cells.Item(27, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 241
$ws.Cells.Item(27, 17).Value = 60
$ws.Cells.Item(27, 18).Value = "Hortaliza"

# Make sure the D27 cell keeps/gets the date-style formatting used by the other date cells
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
